$d = $word.ActiveDocument

function Replace-ParagraphText($paraIndex, $newText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $old = $rng.Text.TrimEnd([char]13, [char]0)
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- Paragraph 1: title/subtitle line ---
# Final run content (per diff): <w:br/><w:t>title</w:t><w:br/><w:t>subtitle</w:t>
$lineBreak = [char]11
$title = '⚡️🚀המאמר היומי של מייק -19.10.24: ⚡️🚀'
$subtitle = 'DiffCSE: Difference-based Contrastive Learning for Sentence Embeddings'
$newP1 = $lineBreak + $title + $lineBreak + $subtitle
Replace-ParagraphText 1 $newP1

# --- Paragraph 2 ---
Replace-ParagraphText 2 'סקירה קצרה ואחרונה(כנראה) במיני-סדרה על איך לבנות ייצוג דאטה באמצעות שיטות למידה ניגודית. כבר הסברתי על הלמידה הניגודית בשתי בסקירות הקודמות. בקצרה, מאמנים מודל הבונה אמבדינג לדאטה המקרב ייצוגים של פיסות דאטה דומות ולהרחיק פיסות דאטה לא דומות. וכאמור הוצעו עשרות שיטות לעשות זאת לדאטה מדומיינים שונים.'

# --- Paragraph 3 ---
Replace-ParagraphText 3 'המאמר מציע שיטת CL העושה זאת בצורה מתוחכמת יותר (לטעמי). הרי אחת המטרות של בניית ייצוג הדאטה היא שהוא ישקף את התכונות האינהרנטיות של הדאטה והמחברים הציעו דרך ״לאכוף״ את זה על הייצוג. הם מאמנים מודל לבניית ייצוג טקסט כך שהמודל ״יבדיל בין מה אמור ומה לא אמור להיות בתוך הטקסט״.'

# --- Paragraph 4 ---
Replace-ParagraphText 4 'איך הם עשו זאת? הם מיסכו כמה טוקנים בטקסט, ביקשו ממודל אחר לחזות את הטוקן הזה ואז אימנו את ייצוג כך שבעזרתו יהיה ניתן להבדיל בין הטוקנים שנחזו ואלו שלא. כלומר בנוסף למודל החיזוי (לא אומן) ומודל לבניית אמבדינג הם אימנו עוד מודל לסיווג בינארי שמטרתו להגיד האם טוקן נחזה או לא. וייצוג הטקסט מוזן למודל הסיווג הזה.'

# --- Paragraph 5 ---
Replace-ParagraphText 5 'דרך אגב פונקציית הלוס למודל הסיווג דומה לזו של GAN אבל אין באמת קשר בין שני הדברים (זה טיפה בלבל אותי בהתחלה)....'

# --- Paragraph 6: becomes the final URL paragraph ---
Replace-ParagraphText 6 'https://arxiv.org/pdf/2204.10298'

# --- Paragraphs 7-10 (original): removed entirely ---
# Delete from the end backwards so indices stay valid.
$d.Paragraphs(10).Range.Delete()
$d.Paragraphs(9).Range.Delete()
$d.Paragraphs(8).Range.Delete()
$d.Paragraphs(7).Range.Delete()
